# Update cryptocurrency price (D) and 1h volume change (E) columns
# with the latest scraped values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.665.54"
$ws.Range("E2").Value = "  +1.15%  "
$ws.Range("D3").Value = "1.564.51"
$ws.Range("E3").Value = "  -0.25%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  -0.40%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "210.32"
$ws.Range("E5").Value = "  -0.27%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.489"
$ws.Range("E6").Value = "  -0.49%  "
$ws.Range("E7").Value = "  -0.46%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "25.11"
$ws.Range("E8").Value = "  +5.55%  "
$ws.Range("E9").Value = "  -0.28%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0586"
$ws.Range("E10").Value = "  -0.09%  "
$ws.Range("E11").Value = "  -0.08%  "
$ws.Range("D12").Value = "1.786.09"
$ws.Range("E12").Value = "  -0.45%  "
$ws.Range("D13").Value = "1.562.78"
$ws.Range("E13").Value = "  -0.54%  "
$ws.Range("D14").Value = "28.677.52"
$ws.Range("E14").Value = "  +1.17%  "
$ws.Range("E15").Value = "  +0.74%  "
$ws.Range("E16").Value = "  -0.69%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.29"
$ws.Range("E17").Value = "  +0.35%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "228.57"
$ws.Range("E18").Value = "  +0.62%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.35"
$ws.Range("E19").Value = "  -0.30%  "
$ws.Range("E20").Value = "  -0.25%  "
$ws.Range("E22").Value = "  -0.32%  "
$ws.Range("E23").Value = "  +0.93%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.08"
$ws.Range("E24").Value = "  +1.13%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "151.33"
$ws.Range("E25").Value = "  +0.38%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "14.76"
$ws.Range("E26").Value = "  -0.73%  "
$ws.Range("E27").Value = "  +0.68%  "
$ws.Range("E28").Value = "  -0.41%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.23"
$ws.Range("E29").Value = "  -1.65%  "
$ws.Range("E30").Value = "  -4.10%  "
$ws.Range("E31").Value = "  -2.56%  "
$ws.Range("E32").Value = "  +0.12%  "
$ws.Range("D33").Value = "1.393.20"
$ws.Range("E33").Value = "  +1.13%  "
$ws.Range("E34").Value = "  -3.21%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.02"
$ws.Range("E35").Value = "  -4.40%  "
$ws.Range("E36").Value = "  -1.26%  "
$ws.Range("E37").Value = "  +1.57%  "
$ws.Range("E38").Value = "  -2.39%  "
$ws.Range("E39").Value = "  -0.55%  "
$ws.Range("E40").Value = "  +1.75%  "
$ws.Range("E41").Value = "  -0.35%  "
$ws.Range("E42").Value = "  -0.37%  "
$ws.Range("E43").Value = "  -1.14%  "
$ws.Range("E44").Value = "  -2.93%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "63.99"
$ws.Range("E45").Value = "  +2.96%  "
$ws.Range("E46").Value = "  -1.77%  "
$ws.Range("D47").Value = "1.698.94"
$ws.Range("E47").Value = "  -0.44%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.870"
$ws.Range("E48").Value = "  -5.04%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "85.10"
$ws.Range("E49").Value = "  -0.33%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "43.22"
$ws.Range("E50").Value = "  +7.31%  "
$ws.Range("E51").Value = "  +0.37%  "
